$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "15 марта"
$ws.Range("B20").Value = "Поиск инструмента для отслеживания утечек памяти, а также написание про него текста в теории вкр, оптимизация копирований при приеме данных"
$ws.Range("B20").HorizontalAlignment = -4152

$ws.Range("B23").Select()
